# Auto-generated script applying Zodiark_Profits market-data refresh to before.xlsx
# Regenerated cell values per sheet, matching the scheduled runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H9").Value = 221.9375
$ws.Range("I9").Value = 220.67857
$ws.Range("K9").Value = 220.67857
$ws.Range("M9").Value = -51.67857000000001
$ws.Range("H15").Value = 1808.2609
$ws.Range("I15").Value = 1808.2609
$ws.Range("K15").Value = 5424.7827
$ws.Range("M15").Value = -5255.7827
$ws.Range("H17").Value = 2770.6
$ws.Range("I17").Value = 1059.5
$ws.Range("J17").Value = 2960.7222
$ws.Range("K17").Value = 3178.5
$ws.Range("L17").Value = 8882.1666
$ws.Range("M17").Value = -3010.5
$ws.Range("N17").Value = -9218.1666
$ws.Range("H32").Value = 6082.5
$ws.Range("J32").Value = 6623.75
$ws.Range("L32").Value = 6623.75
$ws.Range("N32").Value = -7275.75
$ws.Range("H43").Value = 2754.6667
$ws.Range("I43").Value = 3157.6667
$ws.Range("J43").Value = 2620.3333
$ws.Range("K43").Value = 3157.6667
$ws.Range("L43").Value = 2620.3333
$ws.Range("M43").Value = -3088.6667
$ws.Range("N43").Value = -2758.3333
$ws.Range("H51").Value = 4477.636
$ws.Range("I51").Value = 4253
$ws.Range("K51").Value = 4253
$ws.Range("M51").Value = -3769
$ws.Range("H55").Value = 203.44444
$ws.Range("I55").Value = 205.72728
$ws.Range("K55").Value = 205.72728
$ws.Range("M55").Value = 8.272719999999993
$ws.Range("H70").Value = 5268.5713
$ws.Range("I70").Value = 7993
$ws.Range("K70").Value = 23979
$ws.Range("M70").Value = -23709
$ws.Range("H73").Value = 5268.5713
$ws.Range("I73").Value = 7993
$ws.Range("K73").Value = 23979
$ws.Range("M73").Value = -23043
$ws.Range("H76").Value = 7200
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 7200
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H100").Value = 3649.3333
$ws.Range("I100").Value = 2268.9285
$ws.Range("K100").Value = 2268.9285
$ws.Range("M100").Value = -1727.9285
$ws.Range("H113").Value = 8535.191999999999
$ws.Range("J113").Value = 9563.454
$ws.Range("L113").Value = 9563.454
$ws.Range("N113").Value = -16071.454
$ws.Range("H125").Value = 1998
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H135").Value = 693.8946999999999
$ws.Range("I135").Value = 588.8889
$ws.Range("K135").Value = 5300.0001
$ws.Range("M135").Value = -2765.0001
$ws.Range("H137").Value = 7249309.5
$ws.Range("I137").Value = 19236582
$ws.Range("J137").Value = 1191.2325
$ws.Range("K137").Value = 57709746
$ws.Range("L137").Value = 3573.6975
$ws.Range("M137").Value = -57707196
$ws.Range("N137").Value = -8673.6975
$ws.Range("H138").Value = 3881.2046
$ws.Range("J138").Value = 4725.619
$ws.Range("L138").Value = 14176.857
$ws.Range("N138").Value = -24456.857
$ws.Range("H141").Value = 7494.853
$ws.Range("I141").Value = 3642.8667
$ws.Range("J141").Value = 36384.75
$ws.Range("K141").Value = 10928.6001
$ws.Range("L141").Value = 109154.25
$ws.Range("M141").Value = -5748.6001
$ws.Range("N141").Value = -119514.25
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 1502.8572
$ws.Range("I2").Value = 1346.6666
$ws.Range("K2").Value = 1346.6666
$ws.Range("M2").Value = -1233.6666
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -34
$ws.Range("H6").Value = 502
$ws.Range("I6").Value = 502
$ws.Range("K6").Value = 502
$ws.Range("M6").Value = -329
$ws.Range("H32").Value = 3863.047
$ws.Range("I32").Value = 3303.1829
$ws.Range("J32").Value = 19166
$ws.Range("K32").Value = 3303.1829
$ws.Range("L32").Value = 19166
$ws.Range("M32").Value = -3016.1829
$ws.Range("N32").Value = -19740
$ws.Range("H34").Value = 41445
$ws.Range("I34").Value = 41445
$ws.Range("K34").Value = 41445
$ws.Range("M34").Value = -41174
$ws.Range("H40").Value = 49242.5
$ws.Range("I40").Value = 49485
$ws.Range("J40").Value = 49161.668
$ws.Range("K40").Value = 49485
$ws.Range("L40").Value = 49161.668
$ws.Range("M40").Value = -49309
$ws.Range("N40").Value = -49513.668
$ws.Range("H45").Value = 1639.3182
$ws.Range("J45").Value = 1451.2142
$ws.Range("L45").Value = 1451.2142
$ws.Range("N45").Value = -2205.2142
$ws.Range("H61").Value = 3809.75
$ws.Range("I61").Value = 3312
$ws.Range("K61").Value = 3312
$ws.Range("M61").Value = -3100
$ws.Range("H74").Value = 11183.617
$ws.Range("I74").Value = 7715.2964
$ws.Range("K74").Value = 7715.2964
$ws.Range("M74").Value = -6841.2964
$ws.Range("H76").Value = 200168260
$ws.Range("J76").Value = 210319.5
$ws.Range("L76").Value = 210319.5
$ws.Range("N76").Value = -210995.5
$ws.Range("H77").Value = 11183.617
$ws.Range("I77").Value = 7715.2964
$ws.Range("K77").Value = 38576.482
$ws.Range("M77").Value = -34208.482
$ws.Range("H79").Value = 200168260
$ws.Range("J79").Value = 210319.5
$ws.Range("L79").Value = 210319.5
$ws.Range("N79").Value = -212659.5
$ws.Range("H110").Value = 3400.0833
$ws.Range("I110").Value = 2224.1667
$ws.Range("J110").Value = 4576
$ws.Range("K110").Value = 2224.1667
$ws.Range("L110").Value = 4576
$ws.Range("M110").Value = -179.1667000000002
$ws.Range("N110").Value = -8666
$ws.Range("H116").Value = 1502.8572
$ws.Range("I116").Value = 1346.6666
$ws.Range("K116").Value = 1346.6666
$ws.Range("M116").Value = 947.3334
$ws.Range("H132").Value = 11042.743
$ws.Range("I132").Value = 8503.423000000001
$ws.Range("J132").Value = 18378.555
$ws.Range("K132").Value = 25510.269
$ws.Range("L132").Value = 55135.665
$ws.Range("M132").Value = -22980.269
$ws.Range("N132").Value = -60195.665
$ws.Range("H136").Value = 3809.75
$ws.Range("I136").Value = 3312
$ws.Range("K136").Value = 9936
$ws.Range("M136").Value = -7386

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 1502.8572
$ws.Range("I3").Value = 1346.6666
$ws.Range("K3").Value = 1346.6666
$ws.Range("M3").Value = -1232.6666
$ws.Range("H11").Value = 724.5
$ws.Range("J11").Value = 750
$ws.Range("L11").Value = 750
$ws.Range("N11").Value = -1030
$ws.Range("H105").Value = 2391.6
$ws.Range("I105").Value = 2253.3333
$ws.Range("K105").Value = 2253.3333
$ws.Range("M105").Value = -506.3332999999998
$ws.Range("H107").Value = 3646.35
$ws.Range("I107").Value = 3282.5334
$ws.Range("K107").Value = 3282.5334
$ws.Range("M107").Value = -1362.5334
$ws.Range("H134").Value = 2593.8
$ws.Range("I134").Value = 2593.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7781.400000000001
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H141").Value = 211398
$ws.Range("J141").Value = 211398
$ws.Range("L141").Value = 211398
$ws.Range("N141").Value = -221758

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 1965.2142
$ws.Range("J16").Value = 2313.1667
$ws.Range("L16").Value = 2313.1667
$ws.Range("N16").Value = -2887.1667
$ws.Range("H22").Value = 793.2778
$ws.Range("I22").Value = 152.21428
$ws.Range("J22").Value = 3037
$ws.Range("K22").Value = 152.21428
$ws.Range("L22").Value = 3037
$ws.Range("M22").Value = 197.78572
$ws.Range("N22").Value = -3737
$ws.Range("H31").Value = 490841.38
$ws.Range("I31").Value = 13090.111
$ws.Range("K31").Value = 13090.111
$ws.Range("M31").Value = -12795.111
$ws.Range("H34").Value = 490841.38
$ws.Range("I34").Value = 13090.111
$ws.Range("K34").Value = 13090.111
$ws.Range("M34").Value = -12888.111
$ws.Range("H58").Value = 4816.8335
$ws.Range("I58").Value = 5724.75
$ws.Range("J58").Value = 3001
$ws.Range("K58").Value = 5724.75
$ws.Range("L58").Value = 3001
$ws.Range("M58").Value = -5521.75
$ws.Range("N58").Value = -3407
$ws.Range("H86").Value = 76927790
$ws.Range("I86").Value = 100003100
$ws.Range("K86").Value = 100003100
$ws.Range("M86").Value = -100001977
$ws.Range("H89").Value = 76927790
$ws.Range("I89").Value = 100003100
$ws.Range("K89").Value = 500015500
$ws.Range("M89").Value = -500009884
$ws.Range("H99").Value = 15241632
$ws.Range("I99").Value = 6900138
$ws.Range("J99").Value = 55558850
$ws.Range("K99").Value = 6900138
$ws.Range("L99").Value = 55558850
$ws.Range("M99").Value = -6898640
$ws.Range("N99").Value = -55561846
$ws.Range("H105").Value = 15271.947
$ws.Range("J105").Value = 903.5
$ws.Range("L105").Value = 903.5
$ws.Range("N105").Value = -4397.5
$ws.Range("H107").Value = 1060.4445
$ws.Range("I107").Value = 420.42856
$ws.Range("K107").Value = 420.42856
$ws.Range("M107").Value = 1499.57144
$ws.Range("H113").Value = 1965.2142
$ws.Range("J113").Value = 2313.1667
$ws.Range("L113").Value = 2313.1667
$ws.Range("N113").Value = -6653.1667
$ws.Range("H122").Value = 2234.1304
$ws.Range("J122").Value = 2414.25
$ws.Range("L122").Value = 7242.75
$ws.Range("N122").Value = -12142.75
$ws.Range("H126").Value = 15241632
$ws.Range("I126").Value = 6900138
$ws.Range("J126").Value = 55558850
$ws.Range("K126").Value = 20700414
$ws.Range("L126").Value = 166676550
$ws.Range("M126").Value = -20697944
$ws.Range("N126").Value = -166681490
$ws.Range("H132").Value = 2325.348
$ws.Range("I132").Value = 1800.2222
$ws.Range("K132").Value = 5400.6666
$ws.Range("M132").Value = -2870.6666
$ws.Range("H134").Value = 3526.2222
$ws.Range("I134").Value = 3190.2307
$ws.Range("J134").Value = 4399.8
$ws.Range("K134").Value = 9570.6921
$ws.Range("L134").Value = 13199.4
$ws.Range("M134").Value = -7035.6921
$ws.Range("N134").Value = -18269.4
$ws.Range("H136").Value = 4816.8335
$ws.Range("I136").Value = 5724.75
$ws.Range("J136").Value = 3001
$ws.Range("K136").Value = 17174.25
$ws.Range("L136").Value = 9003
$ws.Range("M136").Value = -14624.25
$ws.Range("N136").Value = -14103

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 769.5
$ws.Range("I5").Value = 769.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2308.5
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H7").Value = 412.85715
$ws.Range("I7").Value = 343.75
$ws.Range("J7").Value = 505
$ws.Range("K7").Value = 1031.25
$ws.Range("L7").Value = 1515
$ws.Range("M7").Value = -919.25
$ws.Range("N7").Value = -1739
$ws.Range("H12").Value = 54.266666
$ws.Range("J12").Value = 38.2
$ws.Range("L12").Value = 114.6
$ws.Range("N12").Value = -460.6
$ws.Range("H14").Value = 316.84616
$ws.Range("I14").Value = 316.84616
$ws.Range("K14").Value = 950.5384799999999
$ws.Range("M14").Value = -777.5384799999999
$ws.Range("H23").Value = 225.22223
$ws.Range("I23").Value = 161.33333
$ws.Range("J23").Value = 353
$ws.Range("K23").Value = 483.99999
$ws.Range("L23").Value = 1059
$ws.Range("M23").Value = -248.99999
$ws.Range("N23").Value = -1529
$ws.Range("H60").Value = 854.2857
$ws.Range("I60").Value = 862.3333
$ws.Range("K60").Value = 2586.9999
$ws.Range("M60").Value = -2335.9999
$ws.Range("H107").Value = 990.375
$ws.Range("J107").Value = 1098.4
$ws.Range("L107").Value = 3295.2
$ws.Range("N107").Value = -7135.200000000001
$ws.Range("H113").Value = 9805760
$ws.Range("I113").Value = 2054.923
$ws.Range("J113").Value = 15874720
$ws.Range("K113").Value = 6164.768999999999
$ws.Range("L113").Value = 47624160
$ws.Range("M113").Value = -3994.768999999999
$ws.Range("N113").Value = -47628500
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("H134").Value = 3315.25
$ws.Range("I134").Value = 1998.1111
$ws.Range("J134").Value = 7266.6665
$ws.Range("K134").Value = 5994.3333
$ws.Range("L134").Value = 21799.9995
$ws.Range("M134").Value = -924.3333000000002
$ws.Range("N134").Value = -31939.9995
$ws.Range("H135").Value = 769.5
$ws.Range("I135").Value = 769.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6925.5
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 19000.143
$ws.Range("I70").Value = 31869.05
$ws.Range("J70").Value = 7301.136
$ws.Range("K70").Value = 31869.05
$ws.Range("L70").Value = 7301.136
$ws.Range("M70").Value = -31599.05
$ws.Range("N70").Value = -7841.136
$ws.Range("H73").Value = 19000.143
$ws.Range("I73").Value = 31869.05
$ws.Range("J73").Value = 7301.136
$ws.Range("K73").Value = 31869.05
$ws.Range("L73").Value = 7301.136
$ws.Range("M73").Value = -30933.05
$ws.Range("N73").Value = -9173.136
$ws.Range("H98").Value = 24771.5
$ws.Range("J98").Value = 24771.5
$ws.Range("L98").Value = 24771.5
$ws.Range("N98").Value = -30761.5
$ws.Range("H102").Value = 3205.879
$ws.Range("I102").Value = 2969.8
$ws.Range("J102").Value = 5566.6665
$ws.Range("K102").Value = 2969.8
$ws.Range("L102").Value = 5566.6665
$ws.Range("M102").Value = -1347.8
$ws.Range("N102").Value = -8810.666499999999
$ws.Range("H113").Value = 3151.8333
$ws.Range("I113").Value = 2907
$ws.Range("K113").Value = 2907
$ws.Range("M113").Value = -737
$ws.Range("H122").Value = 3174.5667
$ws.Range("I122").Value = 2562.5
$ws.Range("K122").Value = 7687.5
$ws.Range("M122").Value = -5237.5
$ws.Range("H126").Value = 4186.5815
$ws.Range("I126").Value = 4220.657
$ws.Range("K126").Value = 12661.971
$ws.Range("M126").Value = -10191.971

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H2").Value = 13500
$ws.Range("J2").Value = 13500
$ws.Range("L2").Value = 13500
$ws.Range("N2").Value = -13724
$ws.Range("H55").Value = 1777
$ws.Range("I55").Value = 916.0714
$ws.Range("J55").Value = 3116.2222
$ws.Range("K55").Value = 916.0714
$ws.Range("L55").Value = 3116.2222
$ws.Range("M55").Value = -743.0714
$ws.Range("N55").Value = -3462.2222
$ws.Range("H61").Value = 9765
$ws.Range("I61").Value = 10897.5
$ws.Range("K61").Value = 10897.5
$ws.Range("M61").Value = -10695.5
$ws.Range("H69").Value = 302165.56
$ws.Range("J69").Value = 302165.56
$ws.Range("L69").Value = 302165.56
$ws.Range("N69").Value = -303787.56
$ws.Range("H72").Value = 302165.56
$ws.Range("J72").Value = 302165.56
$ws.Range("L72").Value = 906496.6799999999
$ws.Range("N72").Value = -914608.6799999999
$ws.Range("H113").Value = 9765
$ws.Range("I113").Value = 10897.5
$ws.Range("K113").Value = 10897.5
$ws.Range("M113").Value = -8727.5
$ws.Range("H132").Value = 6924.913
$ws.Range("I132").Value = 7126.222
$ws.Range("K132").Value = 21378.666
$ws.Range("M132").Value = -18848.666
$ws.Range("H134").Value = 98000
$ws.Range("J134").Value = 98000
$ws.Range("L134").Value = 98000
$ws.Range("N134").Value = -108140
$ws.Range("H136").Value = 4544.6855
$ws.Range("I136").Value = 4659.1377
$ws.Range("K136").Value = 13977.4131
$ws.Range("M136").Value = -11427.4131

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H34").Value = 30813
$ws.Range("I34").Value = 29417.334
$ws.Range("K34").Value = 29417.334
$ws.Range("M34").Value = -29214.334
$ws.Range("H37").Value = 31724.5
$ws.Range("I37").Value = 27299.334
$ws.Range("K37").Value = 27299.334
$ws.Range("M37").Value = -27096.334
$ws.Range("H40").Value = 30263
$ws.Range("I40").Value = 30008
$ws.Range("K40").Value = 30008
$ws.Range("M40").Value = -29859
$ws.Range("H45").Value = 8999.666999999999
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("H54").Value = 49976.047
$ws.Range("I54").Value = 49832.332
$ws.Range("K54").Value = 49832.332
$ws.Range("M54").Value = -49312.332
$ws.Range("H64").Value = 70000
$ws.Range("J64").Value = 70000
$ws.Range("L64").Value = 70000
$ws.Range("N64").Value = -70496
$ws.Range("H67").Value = 70000
$ws.Range("J67").Value = 70000
$ws.Range("L67").Value = 70000
$ws.Range("N67").Value = -71716
$ws.Range("H69").Value = 129900
$ws.Range("J69").Value = 129900
$ws.Range("L69").Value = 129900
$ws.Range("N69").Value = -131398
$ws.Range("H72").Value = 129900
$ws.Range("J72").Value = 129900
$ws.Range("L72").Value = 389700
$ws.Range("N72").Value = -397188
$ws.Range("H81").Value = 3870.8572
$ws.Range("I81").Value = 2761.6667
$ws.Range("J81").Value = 10526
$ws.Range("K81").Value = 5523.3334
$ws.Range("L81").Value = 21052
$ws.Range("M81").Value = -4462.3334
$ws.Range("N81").Value = -23174
$ws.Range("H84").Value = 3870.8572
$ws.Range("I84").Value = 2761.6667
$ws.Range("J84").Value = 10526
$ws.Range("K84").Value = 27616.667
$ws.Range("L84").Value = 105260
$ws.Range("M84").Value = -22312.667
$ws.Range("N84").Value = -115868
$ws.Range("H96").Value = 2318075
$ws.Range("I96").Value = 3707639.5
$ws.Range("J96").Value = 2134
$ws.Range("K96").Value = 3707639.5
$ws.Range("L96").Value = 2134
$ws.Range("M96").Value = -3706266.5
$ws.Range("N96").Value = -4880
$ws.Range("H122").Value = 22096274
$ws.Range("I122").Value = 22730212
$ws.Range("K122").Value = 68190636
$ws.Range("M122").Value = -68188186
$ws.Range("H132").Value = 1893.4263
$ws.Range("I132").Value = 1769.4807
$ws.Range("K132").Value = 5308.4421
$ws.Range("M132").Value = -2778.4421
$ws.Range("H136").Value = 4488.9844
$ws.Range("I136").Value = 4278
$ws.Range("K136").Value = 12834
$ws.Range("M136").Value = -10284
$ws.Range("N45").ClearContents()

